# Applies the commit "Add missing dimension tables to code table spreadsheet,
# handle code value conversions" to DimensionalCodeTables.xlsx.
#
# 1. Swap the tab order of PersonAgeType / PersonAgeRangeType (PersonAgeRangeType
#    now comes first).
# 2. Append three new code-table sheets at the end of the workbook:
#       TreatmentProviderType
#       BehavioralHealthEvaluationType
#       MedicationType
# 3. Tidy up view/selection state left over from editing so the final
#    workbook shows MedicationType as the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Swap PersonAgeType / PersonAgeRangeType -----------------------------
$ageRange = $wb.Worksheets.Item("PersonAgeRangeType")
$age = $wb.Worksheets.Item("PersonAgeType")
$ageRange.Move($age)

# --- 2. Add the three new lookup sheets -------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$treatmentProvider = $wb.Worksheets.Add($null, $lastSheet)
$treatmentProvider.Name = "TreatmentProviderType"
$treatmentProvider.Range("A1").Value = "TreatmentProviderTypeID"
$treatmentProvider.Range("B1").Value = "TreatmentProviderTypeDescription"
for ($i = 1; $i -le 10; $i++) {
    $treatmentProvider.Cells.Item($i + 1, 1).Value = $i
    $treatmentProvider.Cells.Item($i + 1, 2).Value = "Provider $i"
}
$treatmentProvider.Cells.Item(12, 1).Value = 99998
$treatmentProvider.Cells.Item(12, 2).Value = "None"
$treatmentProvider.Cells.Item(13, 1).Value = 99999
$treatmentProvider.Cells.Item(13, 2).Value = "Unknown"
$treatmentProvider.Columns.Item(1).ColumnWidth = 28.140625
$treatmentProvider.Columns.Item(2).ColumnWidth = 38.42578125

$behavioralHealth = $wb.Worksheets.Add($null, $treatmentProvider)
$behavioralHealth.Name = "BehavioralHealthEvaluationType"
$behavioralHealth.Range("A1").Value = "BehavioralHealthEvaluationTypeID"
$behavioralHealth.Range("B1").Value = "BehavioralHealthEvaluationTypeDescription"
for ($i = 1; $i -le 10; $i++) {
    $behavioralHealth.Cells.Item($i + 1, 1).Value = $i
    $behavioralHealth.Cells.Item($i + 1, 2).Value = "Diagnosis Category $i"
}
$behavioralHealth.Cells.Item(12, 1).Value = 99998
$behavioralHealth.Cells.Item(12, 2).Value = "None"
$behavioralHealth.Cells.Item(13, 1).Value = 99999
$behavioralHealth.Cells.Item(13, 2).Value = "Unknown"
$behavioralHealth.Columns.Item(1).ColumnWidth = 34.7109375
$behavioralHealth.Columns.Item(2).ColumnWidth = 53.7109375

$medication = $wb.Worksheets.Add($null, $behavioralHealth)
$medication.Name = "MedicationType"
$medication.Range("A1").Value = "MedicationTypeID"
$medication.Range("B1").Value = "MedicationTypeDescription"
for ($i = 1; $i -le 10; $i++) {
    $medication.Cells.Item($i + 1, 1).Value = $i
    $medication.Cells.Item($i + 1, 2).Value = "Medication $i"
}
$medication.Cells.Item(12, 1).Value = 99998
$medication.Cells.Item(12, 2).Value = "None"
$medication.Cells.Item(13, 1).Value = 99999
$medication.Cells.Item(13, 2).Value = "Unknown"
$medication.Columns.Item(1).ColumnWidth = 22.28515625
$medication.Columns.Item(2).ColumnWidth = 42.5703125

# --- 3. Restore/adjust the leftover selections -------------------------------

# ChargeType no longer carries the "last active tab" marker.
$chargeType = $wb.Worksheets.Item("ChargeType")
$chargeType.Range("A1").Select()

# PopulationType's selection moved from A6 to A4:B5.
$populationType = $wb.Worksheets.Item("PopulationType")
$populationType.Range("A4:B5").Select()

# New sheets keep the selection state left by data entry.
$treatmentProvider.Range("A12:B13").Select()
$behavioralHealth.Range("A2:B13").Select()

# MedicationType ends up the active sheet/tab with cell E14 selected.
$medication.Range("E14").Select()

# Scroll the tab strip so ChargeDispositionType is the first visible tab.
$chargeDisposition = $wb.Worksheets.Item("ChargeDispositionType")
$excel.ActiveWindow.ScrollWorkbookTabs($chargeDisposition.Index - 1)
